$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D, rows 1-4 (new column of data)
$ws.Range("D1").Value = "fart"
$ws.Range("E1").Value = "dumb"
$ws.Range("D2").Value = "agree"
$ws.Range("D3").Value = "nah"
$ws.Range("D4").Value = "lol"
$ws.Range("E2").Value = "yoo"
$ws.Range("E4").Value = "damn"
$ws.Range("E3").Value = "haha"

# New row 5 of data
$ws.Range("A5").Value = "gari"
$ws.Range("B5").Value = "chole na"
$ws.Range("C5").Value = "kiu"
$ws.Range("D5").Value = "ok"
$ws.Range("E5").Value = "lalalalal"

# Move selection to A6 as the final active cell (cursor moved past data entry)
$ws.Range("A6").Select()
